# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 222
$ws1.Range("F10").Value = 46
$ws1.Range("F12").Value = 98
$ws1.Range("F13").Value = 2141
$ws1.Range("F14").Value = 53
$ws1.Range("F15").Value = 28
$ws1.Range("F17").Value = 487
$ws1.Range("F18").Value = 151
$ws1.Range("F19").Value = 76
$ws1.Range("F22").Value = 1626
$ws1.Range("F23").Value = 3810
$ws1.Range("F28").Value = 126
$ws1.Range("F29").Value = 2023
$ws1.Range("F31").Value = 464
$ws1.Range("F32").Value = 76

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 223
$ws4.Range("F10").Value = 46
$ws4.Range("F12").Value = 98
$ws4.Range("F13").Value = 2141
$ws4.Range("F14").Value = 53
$ws4.Range("F16").Value = 28
$ws4.Range("F18").Value = 487
$ws4.Range("F19").Value = 151
$ws4.Range("F20").Value = 76
$ws4.Range("F23").Value = 1626
$ws4.Range("F24").Value = 3810
$ws4.Range("F29").Value = 126
$ws4.Range("F30").Value = 2023
$ws4.Range("F32").Value = 464
$ws4.Range("F33").Value = 76
